# Fruta / hortaliza, semanal
# Update weekly price records for Pomelo (Terminal Hortofrutícola Agro Chillán)
# and append the new week's record as row 14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row cell updates (only the columns that actually change for that row).
$updates = @(
    @{ Row = 2;  D = 44210; M = 70; N = 10000; O = 11000; P = 10357; S = 740 },
    @{ Row = 3;  D = 45155; M = 60; N = 15000; O = 15000; P = 15000; S = 1071 },
    @{ Row = 4;  D = 45142; M = 30; N = 15000; O = 15000; P = 15000; S = 1071 },
    @{ Row = 5;  D = 45142; M = 30; N = 14000; O = 14000; P = 14000; Q = "$/caja 14 kilos granel"; S = 1000 },
    @{ Row = 6;  D = 44253; N = 12000; O = 13000; P = 12667; S = 905 },
    @{ Row = 7;  D = 44216; M = 55; N = 11000; O = 12000; P = 11545; S = 825 },
    @{ Row = 10; D = 45138; M = 50; N = 14000; O = 14000; P = 14000; Q = "$/caja 14 kilos granel"; S = 1000 },
    @{ Row = 11; D = 45140; N = 15000; O = 15000; P = 15000; S = 1071 },
    @{ Row = 12; D = 44172; M = 90; N = 8500; O = 9000; P = 8806; Q = "$/caja 14 kilos empedrada"; S = 629 },
    @{ Row = 13; D = 45152; M = 60; N = 16000; O = 16000; P = 16000; Q = "$/caja 14 kilos empedrada"; S = 1143 }
)

foreach ($u in $updates) {
    $rowNum = $u.Row
    foreach ($col in $u.Keys) {
        if ($col -ne "Row") {
            $ws.Range($col + $rowNum).Value2 = $u[$col]
        }
    }
}

# Row 14 is a brand new record; write its values directly, and give the
# date cell D14 the same number format already used by the other date
# cells in column D (e.g. D13) so it stores as a date, not a bare number.
$newRow = @{
    A = 7
    B = "Terminal Hortofrutícola Agro Chillán"
    C = "Ñuble"
    D = 44181
    E = 16
    F = "Fruta"
    G = 100102
    H = "Cítricos"
    I = 100102006
    J = "Pomelo"
    K = "Start Ruby"
    L = "Primera"
    M = 65
    N = 9000
    O = 10000
    P = 9462
    Q = "$/caja 14 kilos empedrada"
    R = "Región de O'Higgins"
    S = 676
    T = 14
}

foreach ($col in $newRow.Keys) {
    $ws.Range($col + "14").Value2 = $newRow[$col]
}

$ws.Range("D14").NumberFormat = $ws.Range("D13").NumberFormat
